# Auto-generated edit script updating cryptos sheet values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "327.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.44%"
$ws.Range("B3").Value = "OKB"
$ws.Range("C3").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.24%"
$ws.Range("B4").Value = "HuobiToken"
$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.506"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.26%"
$ws.Range("B5").Value = "Cronos"
$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08016"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.66%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.004"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.67%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.302"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.90%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9533"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.52%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.570"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-4.54%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1119"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-4.99%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1871"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.32%"
$ws.Range("B12").Value = "MCDex"
$ws.Range("C12").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "10.66"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "26.21%"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09839"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.01%"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04589"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "9.88%"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1066"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.10%"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001280"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.75%"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04083"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.09%"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005897"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.59%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-6.74%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3473"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.39%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1409"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.64%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2545"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.45%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001260"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.60%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004323"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-4.02%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.16%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003741"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.48%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02565"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-3.38%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.02%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007546"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.79%"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.36%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007592"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "13.41%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002019"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.02%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008848"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.66%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007095"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.30%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.52%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "55.03%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003118"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-8.98%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002098"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.52%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001998"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.52%"
